$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.020168084457777
$ws.Cells.Item(2, 4).Value = 1.03005774031357
$ws.Cells.Item(2, 5).Value = 1.021214168540316
$ws.Cells.Item(2, 6).Value = 1.038601733533827
$ws.Cells.Item(2, 9).Value = 1.030159117905667
$ws.Cells.Item(2, 10).Value = 1.025366720859784
$ws.Cells.Item(2, 11).Value = 1.032869854080603
$ws.Cells.Item(2, 12).Value = 1.024052160148313
$ws.Cells.Item(2, 13).Value = 1.041389300736089
$ws.Cells.Item(2, 14).Value = 1.012541577060937

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.021021835154928
$ws.Cells.Item(3, 4).Value = 1.030722766348326
$ws.Cells.Item(3, 5).Value = 1.021935471194385
$ws.Cells.Item(3, 6).Value = 1.039501934825138
$ws.Cells.Item(3, 9).Value = 1.030297665745663
$ws.Cells.Item(3, 10).Value = 1.025858042191998
$ws.Cells.Item(3, 11).Value = 1.033343540563374
$ws.Cells.Item(3, 12).Value = 1.024580101029847
$ws.Cells.Item(3, 13).Value = 1.042099310868047
$ws.Cells.Item(3, 14).Value = 1.012704375368976

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.021574662964624
$ws.Cells.Item(4, 4).Value = 1.031152986988345
$ws.Cells.Item(4, 5).Value = 1.022402949504523
$ws.Cells.Item(4, 6).Value = 1.040084859429011
$ws.Cells.Item(4, 9).Value = 1.030385484641469
$ws.Cells.Item(4, 10).Value = 1.026175726972121
$ws.Cells.Item(4, 11).Value = 1.033649271402065
$ws.Cells.Item(4, 12).Value = 1.024921793160431
$ws.Cells.Item(4, 13).Value = 1.042558537639016
$ws.Cells.Item(4, 14).Value = 1.012809615821165

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.021807164662133
$ws.Cells.Item(5, 4).Value = 1.03133382706031
$ws.Cells.Item(5, 5).Value = 1.022599654589246
$ws.Cells.Item(5, 6).Value = 1.040330023057917
$ws.Cells.Item(5, 9).Value = 1.030421964511353
$ws.Cells.Item(5, 10).Value = 1.02630922467455
$ws.Cells.Item(5, 11).Value = 1.033777613485691
$ws.Cells.Item(5, 12).Value = 1.025065458174989
$ws.Cells.Item(5, 13).Value = 1.042751547709136
$ws.Cells.Item(5, 14).Value = 1.012853834315049

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.021846208125137
$ws.Cells.Item(6, 4).Value = 1.031364189404659
$ws.Cells.Item(6, 5).Value = 1.022632692580483
$ws.Cells.Item(6, 6).Value = 1.040371193045752
$ws.Cells.Item(6, 9).Value = 1.03042806386146
$ws.Cells.Item(6, 10).Value = 1.026331636149832
$ws.Cells.Item(6, 11).Value = 1.033799151670598
$ws.Cells.Item(6, 12).Value = 1.025089581160133
$ws.Cells.Item(6, 13).Value = 1.042783952035982
$ws.Cells.Item(6, 14).Value = 1.012861257340195

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.021577769300517
$ws.Cells.Item(7, 4).Value = 1.031155403481486
$ws.Cells.Item(7, 5).Value = 1.022405577193343
$ws.Cells.Item(7, 6).Value = 1.04008813491843
$ws.Cells.Item(7, 9).Value = 1.030385973814136
$ws.Cells.Item(7, 10).Value = 1.026177511000364
$ws.Cells.Item(7, 11).Value = 1.033650987052438
$ws.Cells.Item(7, 12).Value = 1.02492371275169
$ws.Cells.Item(7, 13).Value = 1.042561116843096
$ws.Cells.Item(7, 14).Value = 1.012810206768015

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.020456531051556
$ws.Cells.Item(8, 4).Value = 1.03028250771716
$ws.Cells.Item(8, 5).Value = 1.021457780407897
$ws.Cells.Item(8, 6).Value = 1.038905869785749
$ws.Cells.Item(8, 9).Value = 1.030206319345949
$ws.Cells.Item(8, 10).Value = 1.025532812720585
$ws.Cells.Item(8, 11).Value = 1.033030098486875
$ws.Cells.Item(8, 12).Value = 1.024230563093935
$ws.Cells.Item(8, 13).Value = 1.041629291844591
$ws.Cells.Item(8, 14).Value = 1.012596616155067

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.018483847322373
$ws.Cells.Item(9, 4).Value = 1.028743699436059
$ws.Cells.Item(9, 5).Value = 1.019793439096704
$ws.Cells.Item(9, 6).Value = 1.036825972779481
$ws.Cells.Item(9, 9).Value = 1.029875762431631
$ws.Cells.Item(9, 10).Value = 1.024395040670599
$ws.Cells.Item(9, 11).Value = 1.031930131327414
$ws.Cells.Item(9, 12).Value = 1.023009807223281
$ws.Cells.Item(9, 13).Value = 1.039985858479927
$ws.Cells.Item(9, 14).Value = 1.012219487934931

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.017170881394464
$ws.Cells.Item(10, 4).Value = 1.027717492238915
$ws.Cells.Item(10, 5).Value = 1.018687876509351
$ws.Cells.Item(10, 6).Value = 1.035441771141199
$ws.Cells.Item(10, 9).Value = 1.029646042127417
$ws.Cells.Item(10, 10).Value = 1.023635432877991
$ws.Cells.Item(10, 11).Value = 1.031192946610795
$ws.Cells.Item(10, 12).Value = 1.022196491418872
$ws.Cells.Item(10, 13).Value = 1.038889365387852
$ws.Cells.Item(10, 14).Value = 1.011967587181848

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.016602881770817
$ws.Cells.Item(11, 4).Value = 1.02727307569208
$ws.Cells.Item(11, 5).Value = 1.018210125089375
$ws.Cells.Item(11, 6).Value = 1.034842987203126
$ws.Cells.Item(11, 9).Value = 1.029544364067562
$ws.Cells.Item(11, 10).Value = 1.023306269851415
$ws.Cells.Item(11, 11).Value = 1.030872835800496
$ws.Cells.Item(11, 12).Value = 1.021844456299204
$ws.Cells.Item(11, 13).Value = 1.038414385303759
$ws.Cells.Item(11, 14).Value = 1.011858402327658

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.016391981425519
$ws.Cells.Item(12, 4).Value = 1.027107992026718
$ws.Cells.Item(12, 5).Value = 1.01803281353639
$ws.Cells.Item(12, 6).Value = 1.034620661770531
$ws.Cells.Item(12, 9).Value = 1.029506265470103
$ws.Cells.Item(12, 10).Value = 1.023183967984852
$ws.Cells.Item(12, 11).Value = 1.030753797773457
$ws.Cells.Item(12, 12).Value = 1.021713716483252
$ws.Cells.Item(12, 13).Value = 1.038237929394268
$ws.Cells.Item(12, 14).Value = 1.011817830077907

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.016437216622839
$ws.Cells.Item(13, 4).Value = 1.027143403347946
$ws.Cells.Item(13, 5).Value = 1.018070840805458
$ws.Cells.Item(13, 6).Value = 1.034668347248383
$ws.Cells.Item(13, 9).Value = 1.029514452706365
$ws.Cells.Item(13, 10).Value = 1.023210203763509
$ws.Cells.Item(13, 11).Value = 1.030779337904665
$ws.Cells.Item(13, 12).Value = 1.021741759615163
$ws.Cells.Item(13, 13).Value = 1.038275780989339
$ws.Cells.Item(13, 14).Value = 1.011826533686585

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.016585447042383
$ws.Cells.Item(14, 4).Value = 1.027259429969173
$ws.Cells.Item(14, 5).Value = 1.01819546546079
$ws.Cells.Item(14, 6).Value = 1.034824607875845
$ws.Cells.Item(14, 9).Value = 1.029541221569804
$ws.Cells.Item(14, 10).Value = 1.023296161069868
$ws.Cells.Item(14, 11).Value = 1.030862998818501
$ws.Cells.Item(14, 12).Value = 1.021833648855186
$ws.Cells.Item(14, 13).Value = 1.038399799948076
$ws.Cells.Item(14, 14).Value = 1.011855048939544

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.016676787339024
$ws.Cells.Item(15, 4).Value = 1.027330916848736
$ws.Cells.Item(15, 5).Value = 1.018272270300104
$ws.Cells.Item(15, 6).Value = 1.034920897129059
$ws.Cells.Item(15, 9).Value = 1.029557670935627
$ws.Cells.Item(15, 10).Value = 1.023349117459645
$ws.Cells.Item(15, 11).Value = 1.03091452727035
$ws.Cells.Item(15, 12).Value = 1.021890267762341
$ws.Cells.Item(15, 13).Value = 1.038476208570461
$ws.Cells.Item(15, 14).Value = 1.011872616002646

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.017208588802299
$ws.Cells.Item(16, 4).Value = 1.027746985543643
$ws.Cells.Item(16, 5).Value = 1.018719603762832
$ws.Cells.Item(16, 6).Value = 1.035481522924966
$ws.Cells.Item(16, 9).Value = 1.029652743724147
$ws.Cells.Item(16, 10).Value = 1.02365727322247
$ws.Cells.Item(16, 11).Value = 1.031214172359106
$ws.Cells.Item(16, 12).Value = 1.022219857804095
$ws.Cells.Item(16, 13).Value = 1.038920884356162
$ws.Cells.Item(16, 14).Value = 1.011974831135649

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.017542314888514
$ws.Cells.Item(17, 4).Value = 1.028007959067524
$ws.Cells.Item(17, 5).Value = 1.019000463641193
$ws.Cells.Item(17, 6).Value = 1.035833346413315
$ws.Cells.Item(17, 9).Value = 1.029711789913035
$ws.Cells.Item(17, 10).Value = 1.023850505666372
$ws.Cells.Item(17, 11).Value = 1.03140189047118
$ws.Cells.Item(17, 12).Value = 1.022426638236584
$ws.Cells.Item(17, 13).Value = 1.039199767624688
$ws.Cells.Item(17, 14).Value = 1.012038918788916

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.017737022065076
$ws.Cells.Item(18, 4).Value = 1.02816017440391
$ws.Cells.Item(18, 5).Value = 1.019164377401058
$ws.Cells.Item(18, 6).Value = 1.036038615379757
$ws.Cells.Item(18, 9).Value = 1.029746017493803
$ws.Cells.Item(18, 10).Value = 1.023963190890299
$ws.Cells.Item(18, 11).Value = 1.031511295864551
$ws.Cells.Item(18, 12).Value = 1.022547262800541
$ws.Cells.Item(18, 13).Value = 1.039362416949156
$ws.Cells.Item(18, 14).Value = 1.012076289375454

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.017803420686538
$ws.Cells.Item(19, 4).Value = 1.028212074811149
$ws.Cells.Item(19, 5).Value = 1.019220283470639
$ws.Cells.Item(19, 6).Value = 1.036108616284578
$ws.Cells.Item(19, 9).Value = 1.029757652049772
$ws.Cells.Item(19, 10).Value = 1.024001609534529
$ws.Cells.Item(19, 11).Value = 1.031548585410024
$ws.Cells.Item(19, 12).Value = 1.022588394832066
$ws.Cells.Item(19, 13).Value = 1.039417872967739
$ws.Cells.Item(19, 14).Value = 1.012089029959321

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.017506504015008
$ws.Cells.Item(20, 4).Value = 1.027979959697012
$ws.Cells.Item(20, 5).Value = 1.018970320408332
$ws.Cells.Item(20, 6).Value = 1.035795593248736
$ws.Cells.Item(20, 9).Value = 1.02970547684774
$ws.Cells.Item(20, 10).Value = 1.023829776126927
$ws.Cells.Item(20, 11).Value = 1.031381759125344
$ws.Cells.Item(20, 12).Value = 1.022404451295548
$ws.Cells.Item(20, 13).Value = 1.039169847998341
$ws.Cells.Item(20, 14).Value = 1.012032043889556

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.016541794658704
$ws.Cells.Item(21, 4).Value = 1.027225263216112
$ws.Cells.Item(21, 5).Value = 1.018158762545625
$ws.Cells.Item(21, 6).Value = 1.034778590528656
$ws.Cells.Item(21, 9).Value = 1.029533347932797
$ws.Cells.Item(21, 10).Value = 1.023270849771678
$ws.Cells.Item(21, 11).Value = 1.030838366469974
$ws.Cells.Item(21, 12).Value = 1.021806589155015
$ws.Cells.Item(21, 13).Value = 1.038363280192773
$ws.Cells.Item(21, 14).Value = 1.011846652350286

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.015935707168446
$ws.Cells.Item(22, 4).Value = 1.026750712635745
$ws.Cells.Item(22, 5).Value = 1.017649352902694
$ws.Cells.Item(22, 6).Value = 1.034139679724096
$ws.Cells.Item(22, 9).Value = 1.029423209931422
$ws.Cells.Item(22, 10).Value = 1.022919222519706
$ws.Cells.Item(22, 11).Value = 1.03049593573523
$ws.Cells.Item(22, 12).Value = 1.021430815684666
$ws.Cells.Item(22, 13).Value = 1.037856002657064
$ws.Cells.Item(22, 14).Value = 1.011729996193985

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.016256961210334
$ws.Cells.Item(23, 4).Value = 1.027002284365435
$ws.Cells.Item(23, 5).Value = 1.017919319553237
$ws.Cells.Item(23, 6).Value = 1.034478328629931
$ws.Cells.Item(23, 9).Value = 1.029481777290499
$ws.Cells.Item(23, 10).Value = 1.023105646110902
$ws.Cells.Item(23, 11).Value = 1.030677538065206
$ws.Cells.Item(23, 12).Value = 1.021630007926847
$ws.Cells.Item(23, 13).Value = 1.038124934367193
$ws.Cells.Item(23, 14).Value = 1.011791846518882

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.017522685251701
$ws.Cells.Item(24, 4).Value = 1.027992611426385
$ws.Cells.Item(24, 5).Value = 1.01898394055078
$ws.Cells.Item(24, 6).Value = 1.035812652105035
$ws.Cells.Item(24, 9).Value = 1.02970833010861
$ws.Cells.Item(24, 10).Value = 1.023839142987763
$ws.Cells.Item(24, 11).Value = 1.03139085588419
$ws.Cells.Item(24, 12).Value = 1.022414476579039
$ws.Cells.Item(24, 13).Value = 1.039183367447144
$ws.Cells.Item(24, 14).Value = 1.012035150393625

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.018993459451563
$ws.Cells.Item(25, 4).Value = 1.029141584966492
$ws.Cells.Item(25, 5).Value = 1.020223013812886
$ws.Cells.Item(25, 6).Value = 1.037363261440536
$ws.Cells.Item(25, 9).Value = 1.029962870951265
$ws.Cells.Item(25, 10).Value = 1.024689379862011
$ws.Cells.Item(25, 11).Value = 1.032215187844333
$ws.Cells.Item(25, 12).Value = 1.012317071403381
$ws.Cells.Item(25, 13).Value = 1.04041088592868
$ws.Cells.Item(25, 14).Value = 1.012317071403381

